$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.30069999999999
$ws.Range("C8").Value = -11.92579999999999
$ws.Range("A12").Value = -21.87610000000002
$ws.Range("C12").Value = -12.9159
$ws.Range("C14").Value = -12.19189999999999
$ws.Range("C22").Value = -11.14569999999999
